$d = $word.ActiveDocument

# 1) "Oder mit einem Versionsverwaltungssystem wie z.B. Git." paragraph:
#    merge the split runs (proofErr-wrapped "Git" + ".") back into one run.
$d.Content.Find.Execute(
    "Oder mit einem Versionsverwaltungssystem wie z.B. Git.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Oder mit einem Versionsverwaltungssystem wie z.B. Git.", 2)

# 2) "Versions-Chaos (resolved)" paragraph:
#    merge its split runs into one, then turn the old closing signature
#    into three paragraphs: the merged line, a blank line, and the new
#    "Mit freundlichen Gruessen" greeting.
$d.Content.Find.Execute(
    "Versions-Chaos (resolved)",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Versions-Chaos (resolved)", 2)

$greeting = "Mit freundlichen Gr" + [char]0x00FC + [char]0x00DF + "en"

$p = $d.Paragraphs(8)
$rng = $p.Range
$rng.Text = "Versions-Chaos (resolved)" + [char]13 + [char]13 + $greeting
